$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row appended below the existing table (row 18 -> row 19).
# Column A holds a plain text date (matches the existing rows, which are
# stored as text, not real Excel dates), so force a text number format
# before assigning the value - otherwise Excel auto-converts a
# "YYYY-MM-DD"-looking string into a date serial number. ClearFormats()
# afterwards drops that temporary number-format override again, so the
# new cell ends up with the same "no explicit style" look as its
# neighbours (A2:A18).
$dateCell = $ws.Cells.Item(19, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2025-09-03"
$dateCell.ClearFormats()

$ws.Cells.Item(19, 2).Value = 57.9900016784668
$ws.Cells.Item(19, 3).Value = 692.0499877929688
$ws.Cells.Item(19, 4).Value = 326.1000061035156
